$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated PSSM values for B2:K21 (supplemental figures update)
$data = @(
    @(-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132),
    @(-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132,1.20669787909877,-19.32476068708132,-19.32476068708132),
    @(-19.32476068708132,-19.32476068708132,1.978233481071232,-19.32476068708132,3.22275167437594,-19.32476068708132,1.506695788410608,-19.32476068708132,0.9366789023801282,-19.32476068708132),
    @(-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132,2.803150848692987,-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132),
    @(-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132),
    @(2.487313573614913,-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132),
    @(-19.32476068708132,-19.32476068708132,-19.32476068708132,1.999642384169151,-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132),
    @(3.847259757213449,-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132),
    @(-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132,1.744548119712841,-19.32476068708132,2.214860575972767),
    @(-19.32476068708132,-19.32476068708132,-19.32476068708132,2.802445565401993,-19.32476068708132,2.834530785677312,-19.32476068708132,-19.32476068708132,-19.32476068708132,1.93477168858323),
    @(-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132),
    @(-19.32476068708132,-19.32476068708132,-19.32476068708132,2.431119946881467,-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132,1.711359923704214,1.76363226019265),
    @(-19.32476068708132,-19.32476068708132,1.369029964797617,-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132,1.95851591917535),
    @(-19.32476068708132,-19.32476068708132,1.631723950703257,-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132),
    @(-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132,1.936138752364368,-19.32476068708132),
    @(-19.32476068708132,4.321926007678873,1.515713764935815,-19.32476068708132,-19.32476068708132,-19.32476068708132,2.049014121823125,2.060573549433241,2.505867860415738,-19.32476068708132),
    @(-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132,-19.32476068708132,1.991200910131685,2.028395293965656,2.407205041510697,-19.32476068708132),
    @(-19.32476068708132,-19.32476068708132,1.991841813024627,-19.32476068708132,-19.32476068708132,-19.32476068708132,1.628629409709961,1.876763220642003,-19.32476068708132,-19.32476068708132),
    @(-19.32476068708132,-19.32476068708132,1.824250178220723,-19.32476068708132,3.414719371116218,-19.32476068708132,1.666772922694892,1.273492678383294,-19.32476068708132,2.088241682603774),
    @(-19.32476068708132,-19.32476068708132,-19.32476068708132,1.860616270964641,-19.32476068708132,2.55759018372663,1.473661276237438,-19.32476068708132,-19.32476068708132,-19.32476068708132)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 2).Value = $rowVals[$c]
    }
}

Write-Output "Updated B2:K21 with supplemental figure values"